$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the existing row 49, pushing everything down
# and creating blank rows 50 and 51 (inheriting formatting from row 49,
# e.g. the date style on column D).
$ws.Range("A50:A51").EntireRow.Insert()

# Row 49 becomes the updated "June Pearl" / "Primera" record.
$ws.Range("A49").Value = 1
$ws.Range("B49").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C49").Value = "Arica y Parinacota"
$ws.Range("D49").Value = 44595
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100103
$ws.Range("H49").Value = "Frutos de hueso (carozo)"
$ws.Range("I49").Value = 100103006
$ws.Range("J49").Value = "Nectarín"
$ws.Range("K49").Value = "June Pearl"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 300
$ws.Range("N49").Value = 19000
$ws.Range("O49").Value = 20000
$ws.Range("P49").Value = 19500
$ws.Range("Q49").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 1083
$ws.Range("T49").Value = 18

# Row 50 is a new "Venus" / "Segunda" record.
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44595
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103006
$ws.Range("J50").Value = "Nectarín"
$ws.Range("K50").Value = "Venus"
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 270
$ws.Range("N50").Value = 18000
$ws.Range("O50").Value = 20000
$ws.Range("P50").Value = 19000
$ws.Range("Q50").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 1056
$ws.Range("T50").Value = 18

# Row 51 is the original "Artic Pride" / "Segunda" record, moved down.
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44544
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100103
$ws.Range("H51").Value = "Frutos de hueso (carozo)"
$ws.Range("I51").Value = 100103006
$ws.Range("J51").Value = "Nectarín"
$ws.Range("K51").Value = "Artic Pride"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 270
$ws.Range("N51").Value = 18000
$ws.Range("O51").Value = 20000
$ws.Range("P51").Value = 19000
$ws.Range("Q51").Value = "`$/caja 18 kilos granel"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 1056
$ws.Range("T51").Value = 18
